# Update "想去人数" (interest count) figures in column F for the
# 展览 sheet and the corresponding rows in the combined 全部类型 sheet.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (first sheet) — rows keyed by their current F value.
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F3").Value = 162
$ws1.Range("F4").Value = 7728
$ws1.Range("F6").Value = 204
$ws1.Range("F10").Value = 442
$ws1.Range("F11").Value = 163
$ws1.Range("F17").Value = 5584
$ws1.Range("F19").Value = 215
$ws1.Range("F20").Value = 991
$ws1.Range("F21").Value = 228
$ws1.Range("F22").Value = 321

# Sheet "全部类型" (fourth sheet) — same events, different row offsets
# because it merges rows from the other sheets.
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F3").Value = 162
$ws4.Range("F4").Value = 7728
$ws4.Range("F6").Value = 204
$ws4.Range("F10").Value = 442
$ws4.Range("F11").Value = 163
$ws4.Range("F18").Value = 5584
$ws4.Range("F21").Value = 215
$ws4.Range("F22").Value = 991
$ws4.Range("F23").Value = 228
$ws4.Range("F24").Value = 321
